$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.000272967848032444
$ws.Range("D2").Value = 0.6425129170340084
$ws.Range("E2").Value = 0.9981179422718006
$ws.Range("B3").Value = 0.0000000003676266793030402
$ws.Range("C3").Value = 0.01301849549774341
$ws.Range("D3").Value = 0.5650646918323966
$ws.Range("E3").Value = 0.8375274494933351
$ws.Range("B4").Value = 0.000000000005738750001259249
$ws.Range("C4").Value = 0.01176308004029253
$ws.Range("D4").Value = 0.4710639843370272
$ws.Range("E4").Value = 0.6754746205633251
$ws.Range("C5").Value = 0.0000002655036418867233
$ws.Range("D5").Value = 0.02552879416942394
$ws.Range("E5").Value = 0.05205069953698283

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.000321815024269444
$ws.Range("D2").Value = 0.7574896145430087
$ws.Range("E2").Value = 1.176729611678621
$ws.Range("B3").Value = 0.0000000004334129076913899
$ws.Range("C3").Value = 0.01534813522821921
$ws.Range("D3").Value = 0.6661821486545001
$ws.Range("E3").Value = 0.9874016974079234
$ws.Range("B4").Value = 0.000000000006765690480558032
$ws.Range("C4").Value = 0.01386806510706791
$ws.Range("D4").Value = 0.5553601592443344
$ws.Range("E4").Value = 0.796349764182275
$ws.Range("C5").Value = 0.0000003130151099233021
$ws.Range("D5").Value = 0.03009713259059886
$ws.Range("E5").Value = 0.06136509209958291

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = 0.0000125495062176749
$ws.Range("C2").Value = 0.009554437463390963
$ws.Range("D2").Value = 1.1705746034569
$ws.Range("E2").Value = 1.468446836336422
$ws.Range("B3").Value = 0.00008531997183802409
$ws.Range("C3").Value = 0.03452244119359665
$ws.Range("D3").Value = 0.8629479954029199
$ws.Range("E3").Value = 1.05851380799705
$ws.Range("B4").Value = 0.0002529878698520342
$ws.Range("C4").Value = 0.009297202904655476
$ws.Range("D4").Value = 0.7270741065098396
$ws.Range("E4").Value = 0.9219790427795368
$ws.Range("B5").Value = 0.00007949285739058928
$ws.Range("C5").Value = 0.02034304433440102
$ws.Range("D5").Value = 1.058507363696061
$ws.Range("E5").Value = 1.211814485320632

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = 0.0001076099402252142
$ws.Range("C2").Value = 0.007453705115974282
$ws.Range("D2").Value = 1.585703061242156
$ws.Range("E2").Value = 1.666371100707123
$ws.Range("B3").Value = 0.0001148298763755546
$ws.Range("C3").Value = 0.02477726001429727
$ws.Range("D3").Value = 0.871165073528946
$ws.Range("E3").Value = 0.9403475543899196
$ws.Range("B4").Value = 0.0007353957242698485
$ws.Range("C4").Value = 0.006997433578934435
$ws.Range("D4").Value = 0.9332746844457358
$ws.Range("E4").Value = 1.175525920959753
$ws.Range("B5").Value = 0.0003950110178780338
$ws.Range("C5").Value = 0.008880409956795551
$ws.Range("D5").Value = 1.528869999002913
$ws.Range("E5").Value = 1.410775496041512
